# Scheduled runner update: refresh market price/profit figures across the
# leve-profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, WVR) with latest values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 3000
$ws.Range("J21").Value = 2000
$ws.Range("L21").Value = 2000
$ws.Range("N21").Value = -2936
$ws.Range("H23").Value = 3000
$ws.Range("J23").Value = 2000
$ws.Range("L23").Value = 2000
$ws.Range("N23").Value = -2468
$ws.Range("H33").Value = 1594
$ws.Range("I33").Value = 1188
$ws.Range("K33").Value = 1188
$ws.Range("M33").Value = -959
$ws.Range("H43").Value = 1996
$ws.Range("J43").Value = 1989
$ws.Range("L43").Value = 1989
$ws.Range("N43").Value = -2127
$ws.Range("H74").Value = 12137.115
$ws.Range("I74").Value = 14321.277
$ws.Range("J74").Value = 7222.75
$ws.Range("K74").Value = 14321.277
$ws.Range("L74").Value = 7222.75
$ws.Range("M74").Value = -13385.277
$ws.Range("N74").Value = -9094.75
$ws.Range("H77").Value = 12137.115
$ws.Range("I77").Value = 14321.277
$ws.Range("J77").Value = 7222.75
$ws.Range("K77").Value = 71606.38499999999
$ws.Range("L77").Value = 36113.75
$ws.Range("M77").Value = -66926.38499999999
$ws.Range("N77").Value = -45473.75
$ws.Range("H116").Value = 10002.23
$ws.Range("I116").Value = 5119.2856
$ws.Range("J116").Value = 15699
$ws.Range("K116").Value = 5119.2856
$ws.Range("L116").Value = 15699
$ws.Range("M116").Value = -1677.2856
$ws.Range("N116").Value = -22583
$ws.Range("H138").Value = 5048.814
$ws.Range("I138").Value = 6477.5713
$ws.Range("J138").Value = 4771
$ws.Range("K138").Value = 19432.7139
$ws.Range("L138").Value = 14313
$ws.Range("M138").Value = -14292.7139
$ws.Range("N138").Value = -24593

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("H32").Value = 1660.2535
$ws.Range("I32").Value = 1692.5217
$ws.Range("K32").Value = 1692.5217
$ws.Range("M32").Value = -1405.5217
$ws.Range("H45").Value = 25551.277
$ws.Range("I45").Value = 40154.91
$ws.Range("J45").Value = 2602.7144
$ws.Range("K45").Value = 40154.91
$ws.Range("L45").Value = 2602.7144
$ws.Range("M45").Value = -39777.91
$ws.Range("N45").Value = -3356.7144
$ws.Range("H63").Value = 5599.6
$ws.Range("I63").Value = 2749
$ws.Range("K63").Value = 2749
$ws.Range("M63").Value = -2063
$ws.Range("H66").Value = 5599.6
$ws.Range("I66").Value = 2749
$ws.Range("K66").Value = 13745
$ws.Range("M66").Value = -10313
$ws.Range("H132").Value = 2025.5358
$ws.Range("I132").Value = 1886.3265
$ws.Range("K132").Value = 5658.979499999999
$ws.Range("M132").Value = -3128.979499999999
$ws.Range("H133").Value = 105664.664
$ws.Range("J133").Value = 105664.664
$ws.Range("L133").Value = 105664.664
$ws.Range("N133").Value = -110724.664
$ws.Range("H135").Value = 82694.86
$ws.Range("J135").Value = 82694.86
$ws.Range("L135").Value = 82694.86
$ws.Range("N135").Value = -92834.86

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 15156412
$ws.Range("I20").Value = 20839008
$ws.Range("J20").Value = 2823.6667
$ws.Range("K20").Value = 20839008
$ws.Range("L20").Value = 2823.6667
$ws.Range("M20").Value = -20838761
$ws.Range("N20").Value = -3317.6667
$ws.Range("H22").Value = 1146.75
$ws.Range("I22").Value = 1334.8
$ws.Range("K22").Value = 1334.8
$ws.Range("M22").Value = -1161.8
$ws.Range("H94").Value = 51283770
$ws.Range("I94").Value = 66668010
$ws.Range("K94").Value = 66668010
$ws.Range("M94").Value = -66667559
$ws.Range("H99").Value = 5969.5
$ws.Range("I99").Value = 2829.5
$ws.Range("K99").Value = 2829.5
$ws.Range("M99").Value = -1331.5
$ws.Range("H105").Value = 8127217
$ws.Range("I105").Value = 418770.72
$ws.Range("K105").Value = 418770.72
$ws.Range("M105").Value = -417023.72
$ws.Range("H134").Value = 3570.4644
$ws.Range("I134").Value = 3559.2632
$ws.Range("K134").Value = 10677.7896
$ws.Range("M134").Value = -8142.7896

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 4459.6
$ws.Range("I12").Value = 1155
$ws.Range("J12").Value = 6662.6665
$ws.Range("K12").Value = 1155
$ws.Range("L12").Value = 6662.6665
$ws.Range("M12").Value = -985
$ws.Range("N12").Value = -7002.6665
$ws.Range("H31").Value = 12507769
$ws.Range("I31").Value = 8711.333000000001
$ws.Range("J31").Value = 17864508
$ws.Range("K31").Value = 8711.333000000001
$ws.Range("L31").Value = 17864508
$ws.Range("M31").Value = -8416.333000000001
$ws.Range("N31").Value = -17865098
$ws.Range("H34").Value = 12507769
$ws.Range("I34").Value = 8711.333000000001
$ws.Range("J34").Value = 17864508
$ws.Range("K34").Value = 8711.333000000001
$ws.Range("L34").Value = 17864508
$ws.Range("M34").Value = -8509.333000000001
$ws.Range("N34").Value = -17864912
$ws.Range("H58").Value = 1692.0938
$ws.Range("I58").Value = 1509.7826
$ws.Range("K58").Value = 1509.7826
$ws.Range("M58").Value = -1306.7826
$ws.Range("H62").Value = 25004100
$ws.Range("I62").Value = 25004100
$ws.Range("K62").Value = 25004100
$ws.Range("M62").Value = -25003476
$ws.Range("H65").Value = 25004100
$ws.Range("I65").Value = 25004100
$ws.Range("K65").Value = 125020500
$ws.Range("M65").Value = -125017380
$ws.Range("H86").Value = 7555
$ws.Range("I86").Value = 7205.7144
$ws.Range("K86").Value = 7205.7144
$ws.Range("M86").Value = -6082.7144
$ws.Range("H89").Value = 7555
$ws.Range("I89").Value = 7205.7144
$ws.Range("K89").Value = 36028.572
$ws.Range("M89").Value = -30412.572
$ws.Range("H132").Value = 13890665
$ws.Range("I132").Value = 1808.4615
$ws.Range("J132").Value = 30304768
$ws.Range("K132").Value = 5425.3845
$ws.Range("L132").Value = 90914304
$ws.Range("M132").Value = -2895.3845
$ws.Range("N132").Value = -90919364
$ws.Range("H136").Value = 1692.0938
$ws.Range("I136").Value = 1509.7826
$ws.Range("K136").Value = 4529.3478
$ws.Range("M136").Value = -1979.3478

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("H107").Value = 946.7778
$ws.Range("J107").Value = 969.5
$ws.Range("L107").Value = 2908.5
$ws.Range("N107").Value = -6748.5
$ws.Range("H131").Value = 10475.6875
$ws.Range("I131").Value = 19402.5
$ws.Range("J131").Value = 1548.875
$ws.Range("K131").Value = 58207.5
$ws.Range("L131").Value = 4646.625
$ws.Range("M131").Value = -53167.5
$ws.Range("N131").Value = -14726.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").ClearContents()
$ws.Range("H70").Value = 167599.61
$ws.Range("I70").Value = 255612.12
$ws.Range("J70").Value = 26779.6
$ws.Range("K70").Value = 255612.12
$ws.Range("L70").Value = 26779.6
$ws.Range("M70").Value = -255342.12
$ws.Range("N70").Value = -27319.6
$ws.Range("H73").Value = 167599.61
$ws.Range("I73").Value = 255612.12
$ws.Range("J73").Value = 26779.6
$ws.Range("K73").Value = 255612.12
$ws.Range("L73").Value = 26779.6
$ws.Range("M73").Value = -254676.12
$ws.Range("N73").Value = -28651.6
$ws.Range("H80").Value = 32260670
$ws.Range("I80").Value = 166668660
$ws.Range("J80").Value = 2755.88
$ws.Range("K80").Value = 166668660
$ws.Range("L80").Value = 2755.88
$ws.Range("M80").Value = -166667662
$ws.Range("N80").Value = -4751.88
$ws.Range("H83").Value = 32260670
$ws.Range("I83").Value = 166668660
$ws.Range("J83").Value = 2755.88
$ws.Range("K83").Value = 833343300
$ws.Range("L83").Value = 13779.4
$ws.Range("M83").Value = -833338308
$ws.Range("N83").Value = -23763.4
$ws.Range("H103").Value = 52663.332
$ws.Range("J103").Value = 52663.332
$ws.Range("L103").Value = 52663.332
$ws.Range("N103").Value = -55007.332
$ws.Range("H126").Value = 11856.667
$ws.Range("I126").Value = 3815.7144
$ws.Range("K126").Value = 11447.1432
$ws.Range("M126").Value = -8977.143199999999
$ws.Range("H136").Value = 19786.5
$ws.Range("J136").Value = 19786.5
$ws.Range("L136").Value = 59359.5
$ws.Range("N136").Value = -64459.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 1070
$ws.Range("I2").Value = 1070
$ws.Range("K2").Value = 1070
$ws.Range("M2").Value = -958
$ws.Range("H4").Value = 655.1
$ws.Range("I4").Value = 912.5
$ws.Range("K4").Value = 912.5
$ws.Range("M4").Value = -799.5
$ws.Range("H31").Value = 9004.25
$ws.Range("I31").Value = 9339
$ws.Range("J31").Value = 8000
$ws.Range("K31").Value = 9339
$ws.Range("L31").Value = 8000
$ws.Range("M31").Value = -8991
$ws.Range("N31").Value = -8696
$ws.Range("H100").Value = 37037604
$ws.Range("I100").Value = 501.5
$ws.Range("J100").Value = 142857900
$ws.Range("K100").Value = 1003
$ws.Range("L100").Value = 285715800
$ws.Range("M100").Value = -462
$ws.Range("N100").Value = -285716882
